$wb = $excel.ActiveWorkbook

$newId = "4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc"

# ---- "Overview" sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newId + ".md"
$wsOverview.Range("B2").Value = "e2e\" + $newId + ".md"
$wsOverview.Range("G2").Value = "2016-08-30 11:05:02"

foreach ($h in $wsOverview.Hyperlinks) {
  $h.TextToDisplay = "e2e\" + $newId + ".md"
}

# ---- "zh-cn" sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newId + ".md"
$wsZhCn.Range("G2").Value = $newId + ".4ab74e330ad8179519427b2cff08588ae293166e.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 11:04:56"

foreach ($h in $wsZhCn.Hyperlinks) {
  $h.TextToDisplay = $newId + ".md"
}

# ---- "de-de" sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newId + ".md"
$wsDeDe.Range("G2").Value = $newId + ".4ab74e330ad8179519427b2cff08588ae293166e.de-de.xlf"

foreach ($h in $wsDeDe.Hyperlinks) {
  $h.TextToDisplay = $newId + ".md"
}
